$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

function Set-Row($row, $country, $values) {
    if ($country -ne $null) {
        $ws.Cells.Item($row, 1).Value = $country
    }
    for ($i = 0; $i -lt $values.Length; $i++) {
        $ws.Cells.Item($row, 2 + $i).Value = $values[$i]
    }
}

# Update the "last updated" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 22 de Abril de 2020 a las 10:52"

# Row 24 - Austria (numbers updated, country unchanged)
Set-Row 24 $null @(14925, 52, 11328, 3087, 176, 19, 510)

# Rows 36-44: country labels rotate + all stats refreshed
Set-Row 36 "Dinamarca" @(7912, 217, 4700, 2842, 81, 0, 370)
Set-Row 37 "Emiratos Arabes Unidos" @(7755, 0, 1443, 6266, 1, 0, 46)
Set-Row 38 "Bielorrusia" @(7281, 558, 769, 6454, 92, 3, 58)
Set-Row 39 "Noruega" @(7241, 0, 32, 7027, 58, 0, 182)
Set-Row 40 "Indonesia" @(7135, 0, 842, 5677, 0, 0, 616)
Set-Row 41 "Chequia" @(7041, 8, 1800, 5037, 80, 3, 204)
Set-Row 42 "Serbia" @(6890, 0, 977, 5783, 101, 0, 130)
Set-Row 43 "Filipinas" @(6710, 111, 693, 5571, 1, 9, 446)
Set-Row 44 "Australia" @(6647, 2, 4920, 1653, 47, 3, 74)

# Rows 52-55: country labels rotate + all stats refreshed
Set-Row 52 "Banglades" @(3772, 390, 92, 3560, 1, 10, 120)
Set-Row 53 "Luxemburgo" @(3618, 0, 670, 2870, 32, 0, 78)
Set-Row 54 "Egipto" @(3490, 0, 870, 2356, 0, 0, 264)
Set-Row 55 "Sudafrica" @(3465, 0, 1055, 2352, 36, 0, 58)

# Row 64 - Kazajistan (numbers updated only)
Set-Row 64 $null @(2047, 52, 505, 1523, 32, 0, 19)

# Row 68 - Uzbekistan (only D/E updated)
$ws.Cells.Item(68, 4).Value = 372
$ws.Cells.Item(68, 5).Value = 1314

# Row 80 - Camerun (only D/E updated)
$ws.Cells.Item(80, 4).Value = 331
$ws.Cells.Item(80, 5).Value = 789

# Row 84 - Hong Kong (numbers updated only)
Set-Row 84 $null @(1034, 4, 678, 352, 8, 0, 4)

# Row 105 - Estado de Palestina (B, C, E updated)
$ws.Cells.Item(105, 2).Value = 474
$ws.Cells.Item(105, 3).Value = 8
$ws.Cells.Item(105, 5).Value = 399

# Row 158 - Uganda (only D/E updated)
$ws.Cells.Item(158, 4).Value = 41
$ws.Cells.Item(158, 5).Value = 20
